# Updated capital structure database
# Applies the updated values for Egypt - Banks (Regional) dataset rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "2" ---
$ws.Range("D2").Value = 0.2545
$ws.Range("E2").Value = 0.2975
$ws.Range("F2").ClearContents()
$ws.Range("K2").Value = 125.6
$ws.Range("L2").Value = 0.4135660190977939
$ws.Range("M2").Value = 24.81
$ws.Range("N2").Value = 0.07139568345323742
$ws.Range("O2").Value = 0.197531847133758
$ws.Range("P2").Value = 24.81
$ws.Range("Q2").Value = 0.07139568345323742
$ws.Range("R2").Value = 0.197531847133758
$ws.Range("U2").Value = 947.1
$ws.Range("V2").Value = 2.72546762589928
$ws.Range("W2").Value = 0.2758458646616541
$ws.Range("X2").Value = 0.07406199560128188
$ws.Range("Y2").Value = 0.2017838690603722
$ws.Range("Z2").Value = -1.224596774193549
$ws.Range("AB2").Value = 0.06947827410982065
$ws.Range("AC2").Value = -0.06947827410982065
$ws.Range("AD2").Value = 83
$ws.Range("AF2").Value = 83
$ws.Range("AG2").Value = -864.1
$ws.Range("AH2").Value = 0.1927990708478513
$ws.Range("AI2").Value = 0.1163279607568325
$ws.Range("AJ2").Value = 1.672667440960124
$ws.Range("AK2").Value = 3.699058219178082

# --- Row 3: Egyptian Gulf Bank (S.A.E) (CASE:EGBE) ---
$ws.Range("D3").Value = 0.294
$ws.Range("E3").Value = 0.214
$ws.Range("K3").Value = 41
$ws.Range("L3").Value = 0.2671009771986971
$ws.Range("M3").Value = 5.71
$ws.Range("N3").Value = 0.03156440022111664
$ws.Range("O3").Value = 0.1392682926829268
$ws.Range("P3").Value = 5.71
$ws.Range("Q3").Value = 0.03156440022111664
$ws.Range("R3").Value = 0.1392682926829268
$ws.Range("U3").Value = 184.5
$ws.Range("V3").Value = 1.019900497512438
$ws.Range("W3").Value = 0.1541353383458647
$ws.Range("X3").Value = 0.06783433348839681
$ws.Range("Y3").Value = 0.08630100485746785
$ws.Range("Z3").Value = -0.3094758064516129
$ws.Range("AB3").Value = 0.06648093228998375
$ws.Range("AC3").Value = -0.06648093228998375
$ws.Range("AD3").Value = 17.6
$ws.Range("AF3").Value = 17.6
$ws.Range("AG3").Value = -166.9
$ws.Range("AH3").Value = 0.08866498740554157
$ws.Range("AI3").Value = 0.05233422539399346
$ws.Range("AJ3").Value = -11.92142857142857
$ws.Range("AK3").Value = -1.099472990777339

# --- Row 4: alBaraka Bank Egypt S.A.E. (CASE:SAUD) ---
$ws.Range("B4").Value = "alBaraka Bank Egypt S.A.E. (CASE:SAUD)"
$ws.Range("D4").Value = 0.215
$ws.Range("E4").Value = 0.381
$ws.Range("F4").ClearContents()
$ws.Range("K4").Value = 84.59999999999999
$ws.Range("L4").Value = 0.563249001331558
$ws.Range("M4").Value = 19.1
$ws.Range("N4").Value = 0.1146458583433373
$ws.Range("O4").Value = 0.2257683215130024
$ws.Range("P4").Value = 19.1
$ws.Range("Q4").Value = 0.1146458583433373
$ws.Range("R4").Value = 0.2257683215130024
$ws.Range("U4").Value = 762.6
$ws.Range("V4").Value = 4.577430972388956
$ws.Range("W4").Value = 0.3975563909774436
$ws.Range("X4").Value = 0.08028965771416693
$ws.Range("Y4").Value = 0.3172667332632766
$ws.Range("Z4").Value = 0.6056451612903224
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.07247561592965754
$ws.Range("AC4").Value = -0.07247561592965754
$ws.Range("AD4").Value = 65.40000000000001
$ws.Range("AF4").Value = 65.40000000000001
$ws.Range("AG4").Value = -697.2
$ws.Range("AH4").Value = 0.281896551724138
$ws.Range("AI4").Value = 0.1733828207847296
$ws.Range("AJ4").Value = 1.313984168865435
$ws.Range("AK4").Value = 1.809029579657499
